$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 is the "Android  Citizen Data App" task.
# Status (B3) changes from "Proposed" to a brand-new value "Scheduled".
$ws.Range("B3").Value = "Scheduled"

# DueAsDate (D3) now has an actual date value (2019-10-15 => serial 43753).
$ws.Range("D3").Value = "10/15/2019"

# DueDate (C3) becomes a formula deriving the display string from D3, and
# reverts to the default (unstyled) cell format.
$ws.Range("C3").Style = "Normal"
$ws.Range("C3").Formula = '=TEXT(D3,"DD-MMM-YY")'

# Estimate (E3) now has a value.
$ws.Range("E3").Value = 60

# Update the active selection to C7, matching the saved view state.
$ws.Range("C7").Select()
